# Update "想去人数" (interested headcount) figures for several events.
# Sheet "展览" (sheet1) and sheet "全部类型" (sheet4) each list the same
# events; bump the relevant F-column counts on both sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 158
$ws1.Range("F5").Value = 1794
$ws1.Range("F9").Value = 2173
$ws1.Range("F12").Value = 1370
$ws1.Range("F19").Value = 41
$ws1.Range("F22").Value = 21
$ws1.Range("F23").Value = 1181
$ws1.Range("F24").Value = 9

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 158
$ws4.Range("F5").Value = 1794
$ws4.Range("F10").Value = 2173
$ws4.Range("F13").Value = 1370
$ws4.Range("F20").Value = 41
$ws4.Range("F23").Value = 21
$ws4.Range("F24").Value = 1181
$ws4.Range("F25").Value = 9
